# =====================================================================
# [Meeting Minutes][Fix] Next meeting: 11am on Tuesday 3rd Feb;
# Subsequent meetings Mondays biweekly.
# =====================================================================
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1 (Decisions Made -> A: Supervisory meeting date/time):
# "Date/time: 11:00 am on Tuesday." becomes a new, richer sentence and
# a brand-new bullet about subsequent biweekly meetings is appended
# right after it.
# ---------------------------------------------------------------------
$searchRange = $d.Content
$found = $searchRange.Find.Execute("Date/time: 11:00 am on Tuesday.", $true, $false, $false, $false, $false, `
                                    $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find 'Date/time: 11:00 am on Tuesday.'" }

$targetParaStart = $searchRange.Start
$targetPara = $d.Range($targetParaStart, $targetParaStart).Paragraphs(1)

# Insert the brand-new "Subsequent biweekly meetings" bullet straight
# after the target paragraph FIRST -- before we touch / format any text
# inside the target paragraph -- so the new paragraph mark does not
# inherit any of the superscript formatting we are about to add.
$endOfTargetPara = $targetPara.Range.End - 1
$insPt = $d.Range($endOfTargetPara, $endOfTargetPara)
$insPt.InsertParagraphAfter()
$newPara = $d.Range($endOfTargetPara, $endOfTargetPara).Paragraphs(1).Next()
$newParaRange = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$newParaRange.Text = "Subsequent biweekly meetings: 11:00 am on Monday."

# Now rebuild the original bullet's text as a run-by-run sentence.
$textRange = $d.Range($targetParaStart, $targetParaStart)
$textRange = $d.Range($targetParaStart, $targetPara.Range.End - 1)
$textRange.Text = "Upcoming"
$textRange.Collapse(0)
$textRange.InsertAfter(" meeting")
$textRange.Collapse(0)
$textRange.InsertAfter(": 11:00 am on Tuesday")
$textRange.Collapse(0)
$textRange.InsertAfter(", 3")
$textRange.Collapse(0)
$ordStart = $textRange.Start
$textRange.InsertAfter("rd")
$ordRange = $d.Range($ordStart, $ordStart + 2)
$ordRange.Font.Superscript = $true
$tailRange = $d.Range($ordStart + 2, $ordStart + 2)
$tailRange.InsertAfter(" ")
$tailRange.Collapse(0)
$tailRange.InsertAfter("February;")

# ---------------------------------------------------------------------
# Change 2 (Action Items heading):
# Drop the stale "lastRenderedPageBreak" cache marker that sits on the
# "4. " run -- delete + retype the run so no rendering artefact survives.
# ---------------------------------------------------------------------
$searchRange = $d.Content
$found = $searchRange.Find.Execute("4. ", $true, $false, $false, $false, $false, `
                                    $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find the '4. ' action-items run" }
$headingRunStart = $searchRange.Start
$headingRun = $d.Range($headingRunStart, $headingRunStart + 3)
$headingRun.Delete()
$headingInsPt = $d.Range($headingRunStart, $headingRunStart)
$headingInsPt.Font.Bold = $true
$headingInsPt.InsertAfter("4. ")

# ---------------------------------------------------------------------
# Change 3 (Action Items table, Deadline column, both rows):
# "3rd Jan" -> "3rd Feb" in both task rows.
# ---------------------------------------------------------------------
for ($i = 0; $i -lt 2; $i++) {
    $searchRange = $d.Content
    $found = $searchRange.Find.Execute(" Jan", $true, $false, $false, $false, $false, `
                                        $true, 1, $false, "", 0)
    if (-not $found) { throw "Could not find ' Jan' occurrence $i" }
    $searchRange.Text = " "
    $searchRange.Collapse(0)
    $searchRange.InsertAfter("Feb")
}

# ---------------------------------------------------------------------
# Change 4 (5. Next Meeting -> Date):
# "11:00 am, Tuesday 3rd January." becomes
# "11:00 am on Tuesday, 3rd of February 2026." (with "February"/" 2026"
# picking up Arial / #222222 styling, matching a pasted-in fragment).
# ---------------------------------------------------------------------
$searchRange = $d.Content
$found = $searchRange.Find.Execute("11:00 am, Tuesday 3", $true, $false, $false, $false, $false, `
                                    $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find '11:00 am, Tuesday 3'" }
$textRange = $searchRange
$textRange.Text = "11:00 am"
$textRange.Collapse(0)
$textRange.InsertAfter(" on ")
$textRange.Collapse(0)
$textRange.InsertAfter("Tuesday")
$textRange.Collapse(0)
$textRange.InsertAfter(",")
$textRange.Collapse(0)
$textRange.InsertAfter(" 3")
# NOTE: the existing superscript "rd" run right after this point is left
# untouched, matching the diff (it is unchanged context).

$searchRange2 = $d.Content
$found2 = $searchRange2.Find.Execute(" January", $true, $false, $false, $false, $false, `
                                      $true, 1, $false, "", 0)
if (-not $found2) { throw "Could not find ' January'" }
$textRange2 = $searchRange2
$textRange2.Text = " "
$textRange2.Collapse(0)
$textRange2.InsertAfter("of ")
$textRange2.Collapse(0)
$febStart = $textRange2.Start
$textRange2.InsertAfter("February")
$febRange = $d.Range($febStart, $febStart + 8)
$febRange.Font.Name = "Arial"
$febRange.Font.Color = 2236962
$yearInsPt = $d.Range($febStart + 8, $febStart + 8)
$yearStart = $yearInsPt.Start
$yearInsPt.InsertAfter(" 2026")
$yearRange = $d.Range($yearStart, $yearStart + 5)
$yearRange.Font.Name = "Arial"
$yearRange.Font.Color = 2236962

Write-Output "Done."
Write-Output ("Bullet A: " + $targetPara.Range.Text)
Write-Output ("Next meeting date line: " + $d.Paragraphs($d.Paragraphs.Count - 1).Range.Text)
